# Auto-applied update: refresh market-price derived columns (H-N) across Leve profit sheets
# per scheduled runner diff (Sargatanas_Profits.xlsx).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3874.75
$ws.Range("J40").Value = 4857
$ws.Range("L40").Value = 4857
$ws.Range("N40").Value = -5207
$ws.Range("H55").Value = 204.1579
$ws.Range("I55").Value = 210.5
$ws.Range("K55").Value = 210.5
$ws.Range("M55").Value = 3.5
$ws.Range("H98").Value = 83338850
$ws.Range("I98").Value = 83338850
$ws.Range("K98").Value = 83338850
$ws.Range("M98").Value = -83337352
$ws.Range("H111").Value = 7815578.5
$ws.Range("I111").Value = 12502984
$ws.Range("J111").Value = 3234.8333
$ws.Range("K111").Value = 37508952
$ws.Range("L111").Value = 9704.499899999999
$ws.Range("M111").Value = -37505885
$ws.Range("N111").Value = -15838.4999
$ws.Range("H122").Value = 83338850
$ws.Range("I122").Value = 83338850
$ws.Range("K122").Value = 250016550
$ws.Range("M122").Value = -250014100
$ws.Range("H132").Value = 1840.0588
$ws.Range("I132").Value = 1840.0588
$ws.Range("K132").Value = 5520.1764
$ws.Range("M132").Value = -2990.1764
$ws.Range("H138").Value = 8764231
$ws.Range("J138").Value = 12268519
$ws.Range("L138").Value = 36805557
$ws.Range("N138").Value = -36815837

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2276897.5
$ws.Range("I32").Value = 2276897.5
$ws.Range("K32").Value = 2276897.5
$ws.Range("M32").Value = -2276610.5
$ws.Range("H61").Value = 47626270
$ws.Range("I61").Value = 1389.1111
$ws.Range("K61").Value = 1389.1111
$ws.Range("M61").Value = -1177.1111
$ws.Range("H62").Value = 38750
$ws.Range("J62").Value = 38750
$ws.Range("L62").Value = 38750
$ws.Range("N62").Value = -39998
$ws.Range("H65").Value = 38750
$ws.Range("J65").Value = 38750
$ws.Range("L65").Value = 116250
$ws.Range("N65").Value = -122490
$ws.Range("H74").Value = 137584.42
$ws.Range("I74").Value = 401249.75
$ws.Range("J74").Value = 5751.75
$ws.Range("K74").Value = 401249.75
$ws.Range("L74").Value = 5751.75
$ws.Range("M74").Value = -400375.75
$ws.Range("N74").Value = -7499.75
$ws.Range("H77").Value = 137584.42
$ws.Range("I77").Value = 401249.75
$ws.Range("J77").Value = 5751.75
$ws.Range("K77").Value = 2006248.75
$ws.Range("L77").Value = 28758.75
$ws.Range("M77").Value = -2001880.75
$ws.Range("N77").Value = -37494.75
$ws.Range("H94").Value = 35192.832
$ws.Range("J94").Value = 35192.832
$ws.Range("L94").Value = 35192.832
$ws.Range("N94").Value = -36994.832
$ws.Range("H136").Value = 47626270
$ws.Range("I136").Value = 1389.1111
$ws.Range("K136").Value = 4167.3333
$ws.Range("M136").Value = -1617.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 11587.571
$ws.Range("I75").Value = 3222.6
$ws.Range("J75").Value = 32500
$ws.Range("K75").Value = 3222.6
$ws.Range("L75").Value = 32500
$ws.Range("M75").Value = -2286.6
$ws.Range("N75").Value = -34372
$ws.Range("H78").Value = 11587.571
$ws.Range("I78").Value = 3222.6
$ws.Range("J78").Value = 32500
$ws.Range("K78").Value = 9667.799999999999
$ws.Range("L78").Value = 97500
$ws.Range("M78").Value = -4987.799999999999
$ws.Range("N78").Value = -106860
$ws.Range("H107").Value = 28130296
$ws.Range("I107").Value = 34096284
$ws.Range("K107").Value = 34096284
$ws.Range("M107").Value = -34094364
$ws.Range("H134").Value = 5213182
$ws.Range("I134").Value = 8066150
$ws.Range("K134").Value = 24198450
$ws.Range("M134").Value = -24195915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1003.1667
$ws.Range("I22").Value = 881.7778
$ws.Range("J22").Value = 1367.3334
$ws.Range("K22").Value = 881.7778
$ws.Range("L22").Value = 1367.3334
$ws.Range("M22").Value = -531.7778
$ws.Range("N22").Value = -2067.3334
$ws.Range("H31").Value = 5882.193
$ws.Range("I31").Value = 2289.56
$ws.Range("J31").Value = 8688.9375
$ws.Range("K31").Value = 2289.56
$ws.Range("L31").Value = 8688.9375
$ws.Range("M31").Value = -1994.56
$ws.Range("N31").Value = -9278.9375
$ws.Range("H34").Value = 5882.193
$ws.Range("I34").Value = 2289.56
$ws.Range("J34").Value = 8688.9375
$ws.Range("K34").Value = 2289.56
$ws.Range("L34").Value = 8688.9375
$ws.Range("M34").Value = -2087.56
$ws.Range("N34").Value = -9092.9375
$ws.Range("H58").Value = 9617
$ws.Range("I58").Value = 1772
$ws.Range("J58").Value = 11360.333
$ws.Range("K58").Value = 1772
$ws.Range("L58").Value = 11360.333
$ws.Range("M58").Value = -1569
$ws.Range("N58").Value = -11766.333
$ws.Range("H134").Value = 5048.6123
$ws.Range("I134").Value = 2454.12
$ws.Range("J134").Value = 7751.2085
$ws.Range("K134").Value = 7362.36
$ws.Range("L134").Value = 23253.6255
$ws.Range("M134").Value = -4827.36
$ws.Range("N134").Value = -28323.6255
$ws.Range("H136").Value = 9617
$ws.Range("I136").Value = 1772
$ws.Range("J136").Value = 11360.333
$ws.Range("K136").Value = 5316
$ws.Range("L136").Value = 34080.999
$ws.Range("M136").Value = -2766
$ws.Range("N136").Value = -39180.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 585.1429000000001
$ws.Range("I18").Value = 349.33334
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 1048.00002
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -879.0000199999999
$ws.Range("N18").Value = -6338
$ws.Range("H68").Value = 20003672
$ws.Range("I68").Value = 40000940
$ws.Range("K68").Value = 120002820
$ws.Range("M68").Value = -120002009
$ws.Range("H71").Value = 20003672
$ws.Range("I71").Value = 40000940
$ws.Range("K71").Value = 360008460
$ws.Range("M71").Value = -360004404
$ws.Range("H107").Value = 10527082
$ws.Range("I107").Value = 409.6
$ws.Range("J107").Value = 22223384
$ws.Range("K107").Value = 1228.8
$ws.Range("L107").Value = 66670152
$ws.Range("M107").Value = 691.1999999999998
$ws.Range("N107").Value = -66673992
$ws.Range("H139").Value = 80578
$ws.Range("I139").Value = 86459.5
$ws.Range("K139").Value = 259378.5
$ws.Range("M139").Value = -254238.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4163.5
$ws.Range("I7").Value = 2618.875
$ws.Range("K7").Value = 2618.875
$ws.Range("M7").Value = -2506.875
$ws.Range("H122").Value = 3565.7742
$ws.Range("I122").Value = 2612.5715
$ws.Range("J122").Value = 5567.5
$ws.Range("K122").Value = 7837.7145
$ws.Range("L122").Value = 16702.5
$ws.Range("M122").Value = -5387.7145
$ws.Range("N122").Value = -21602.5
$ws.Range("H126").Value = 4163.5
$ws.Range("I126").Value = 2618.875
$ws.Range("K126").Value = 7856.625
$ws.Range("M126").Value = -5386.625
